# "rerun dist commute with harmonised education"
# Updates the dist_commute_r2_ML (D) / dist_commute_r2_LR (E) columns on the
# "accuracy_all" sheet with re-run values, adds a new "All" summary row (23)
# to both sheets, applies a 3-decimal number format to the recomputed D/E
# columns on the "accuracy_all (2)" sheet, and refreshes the sheet selections.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("accuracy_all")
$ws2 = $wb.Worksheets.Item("accuracy_all (2)")

# ---------------------------------------------------------------------------
# 1. Re-run results: updated dist_commute_r2 columns (D = ML, E = LR) on the
#    source sheet. accuracy_all (2) references these cells via IF(...) formulas
#    so its values recalculate automatically.
# ---------------------------------------------------------------------------
$ws1.Range("D2").Value  = 0.17899999999999999
$ws1.Range("D3").Value  = -0.007
$ws1.Range("E3").Value  = -0.01
$ws1.Range("E4").Value  = 0.029
$ws1.Range("D5").Value  = -0.035
$ws1.Range("E5").Value  = 0.013
$ws1.Range("D6").Value  = -0.044
$ws1.Range("E6").Value  = -0.011
$ws1.Range("D7").Value  = 0.06
$ws1.Range("E7").Value  = 0.077
$ws1.Range("D8").Value  = 0.009
$ws1.Range("E8").Value  = 0.039
$ws1.Range("D9").Value  = 0.068
$ws1.Range("E9").Value  = -0.013
$ws1.Range("D10").Value = 0.223
$ws1.Range("E10").Value = 0.253
$ws1.Range("D11").Value = 0.13600000000000001
$ws1.Range("E11").Value = 0.20200000000000001
$ws1.Range("D12").Value = 0.11700000000000001
$ws1.Range("D13").Value = 0.157
$ws1.Range("E13").Value = 0.151
$ws1.Range("D14").Value = 0.216
$ws1.Range("E14").Value = 0.20499999999999999
$ws1.Range("D15").Value = 0.193
$ws1.Range("E15").Value = 0.17299999999999999
$ws1.Range("E16").Value = 0.191
$ws1.Range("D18").Value = 0.17699999999999999
$ws1.Range("E18").Value = 0.154
$ws1.Range("D19").Value = 0.17299999999999999
$ws1.Range("E19").Value = 0.13800000000000001
$ws1.Range("D20").Value = 0.115
$ws1.Range("E20").Value = 0.13400000000000001
$ws1.Range("D21").Value = 0.219
$ws1.Range("D22").Value = 0.11

# ---------------------------------------------------------------------------
# 2. New "All" summary row (row 23) on accuracy_all, plain computed values.
# ---------------------------------------------------------------------------
$ws1.Range("A23").Value = "All"
$ws1.Range("B23").Value = 0.59099999999999997
$ws1.Range("C23").Value = 0.56000000000000005
$ws1.Range("D23").Value = 0.18099999999999999
$ws1.Range("E23").Value = 0.13800000000000001
$ws1.Range("F23").Value = 0.75800000000000001
$ws1.Range("G23").Value = 0.67800000000000005
$ws1.Range("H23").Value = 0.78600000000000003
$ws1.Range("I23").Value = 0.76200000000000001

# ---------------------------------------------------------------------------
# 3. Matching "All" row on accuracy_all (2): A/B/C/F/G/H/I copied as plain
#    values, D/E kept as the usual cross-sheet IF(...) formulas.
# ---------------------------------------------------------------------------
$ws2.Range("A23").Value = "All"
$ws2.Range("B23").Value = 0.59099999999999997
$ws2.Range("C23").Value = 0.56000000000000005
$ws2.Range("D23").Formula = "=IF(accuracy_all!D23>0,accuracy_all!D23,""<0"")"
$ws2.Range("E23").Formula = "=IF(accuracy_all!E23>0,accuracy_all!E23,""<0"")"
$ws2.Range("F23").Value = 0.75800000000000001
$ws2.Range("G23").Value = 0.67800000000000005
$ws2.Range("H23").Value = 0.78600000000000003
$ws2.Range("I23").Value = 0.76200000000000001

# ---------------------------------------------------------------------------
# 4. Apply the new 3-decimal number format to the recomputed dist_commute_r2
#    columns (D,E) for all data rows, including the new "All" row.
# ---------------------------------------------------------------------------
$ws2.Range("D2:E23").NumberFormat = "0.000"

# ---------------------------------------------------------------------------
# 5. Refresh sheet selections / active sheet to match the saved view state.
# ---------------------------------------------------------------------------
$ws1.Range("A23:I23").Select()
$ws2.Activate()
$ws2.Range("D20").Select()
